# Financials update: apply the yearly-financials corrections to the CVVUF sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CVVUF")

# Income Statement section
$ws.Range("J12").Value = 3600
$ws.Range("I14").Value = 200
$ws.Range("E17").Value = 2100
$ws.Range("G17").Value = 1000
$ws.Range("J17").Value = 5200
$ws.Range("E18").Value = -2100
$ws.Range("G18").Value = -1000
$ws.Range("J18").Value = -5200
$ws.Range("G20").Value = 1600
$ws.Range("D21").Value = -800
$ws.Range("F23").Value = -1200
$ws.Range("J23").Value = -5100
$ws.Range("F26").Value = -1200
$ws.Range("J26").Value = -5100
$ws.Range("F27").Value = -1200
$ws.Range("J27").Value = -5100
$ws.Range("G32").Value = -1600
$ws.Range("F33").Value = -1200
$ws.Range("J33").Value = -5100
$ws.Range("F35").Value = -1200
$ws.Range("J35").Value = -5100

# Balance Sheet section
$ws.Range("I41").Value = 900
$ws.Range("J41").Value = 3300
$ws.Range("E42").Value = 400
$ws.Range("D43").Value = 0
$ws.Range("H46").Value = 1100
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 3600
$ws.Range("H48").Value = 800
$ws.Range("D54").Value = 1300
$ws.Range("E54").Value = 1500
$ws.Range("F54").Value = 1700
$ws.Range("G54").Value = 2600
$ws.Range("H54").Value = 2100
$ws.Range("I54").Value = 2400
$ws.Range("J54").Value = 5300
$ws.Range("E57").Value = 100
$ws.Range("J57").Value = 1300
$ws.Range("E60").Value = 100
$ws.Range("J60").Value = 1300
$ws.Range("E66").Value = 100
$ws.Range("J66").Value = 1300
$ws.Range("D72").Value = -55000
$ws.Range("E72").Value = -54200
$ws.Range("F72").Value = -53200
$ws.Range("G72").Value = -60100
$ws.Range("H72").Value = -60700
$ws.Range("I72").Value = -52200
$ws.Range("J72").Value = -50600
$ws.Range("G76").Value = 2500
$ws.Range("H76").Value = 1800
$ws.Range("J76").Value = 3900

# Cash Flow Statement section
$ws.Range("F81").Value = -1200
$ws.Range("J81").Value = -5100
$ws.Range("E89").Value = -1100
$ws.Range("G89").Value = -900
$ws.Range("I89").Value = -2400
$ws.Range("J89").Value = -4900
$ws.Range("G94").Value = 1300
$ws.Range("I102").Value = -2300
$ws.Range("J102").Value = -3900
